$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 38, pushing the existing rows 38..109 down to 39..110
# (matches the "A1:R109" -> "A1:R110" dimension growth in the diff).
$ws.Rows(38).Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Cells.Item(38, 1).Value = 10
$ws.Cells.Item(38, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value = "La Araucanía"
$ws.Cells.Item(38, 4).Value = 45281
$ws.Cells.Item(38, 5).Value = 9
$ws.Cells.Item(38, 6).Value = 100112026
$ws.Cells.Item(38, 7).Value = "Haba"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 180
$ws.Cells.Item(38, 11).Value = 12000
$ws.Cells.Item(38, 12).Value = 13000
$ws.Cells.Item(38, 13).Value = 12500
$ws.Cells.Item(38, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(38, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(38, 16).Value = 500
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"
